# "alteração nas primeiras funções"
#
# Edits applied to the first rows of Plan1:
#  1) The review text in A2 gets extra inner whitespace: "alguns elementos"
#     becomes "alguns" + 18 spaces + "elementos".
#  2) A2's cell format is normalized to the same word-wrap alignment already
#     used by the rest of column A (e.g. A3) - this drops the stray
#     "apply fill" that the previous format carried.
#  3) The active selection moves from B3 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A2")

# 1) Text update (use Value2 - Value returns the COM property descriptor
#    rather than invoking the getter in this host).
$text = $cell.Value2
$cell.Value2 = $text -replace "alguns elementos", "alguns                  elementos"

# 2) Re-apply the plain wrap-text alignment (no fill), matching column A's
#    other cells.
$ws.Range("A3").Copy()
$cell.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Move the selection to A2.
$null = $cell.Select()
